# Scheduled-runner market data refresh: updates currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ columns
# (H:N) for the rows whose Universalis price snapshot changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 337.15384
$ws.Range("J28").Value = 966.3333
$ws.Range("L28").Value = 966.3333
$ws.Range("N28").Value = -1936.3333

$ws.Range("H98").Value = 726.4666999999999
$ws.Range("I98").Value = 766.4167
$ws.Range("J98").Value = 566.6667
$ws.Range("K98").Value = 766.4167
$ws.Range("L98").Value = 566.6667
$ws.Range("M98").Value = 731.5833
$ws.Range("N98").Value = -3562.6667

$ws.Range("H107").Value = 403.82352
$ws.Range("I107").Value = 357.85715
$ws.Range("K107").Value = 357.85715
$ws.Range("M107").Value = 1562.14285

$ws.Range("H112").Value = 1148.4615
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 1134.4
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 3403.2
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -5619.200000000001

$ws.Range("H122").Value = 726.4666999999999
$ws.Range("I122").Value = 766.4167
$ws.Range("J122").Value = 566.6667
$ws.Range("K122").Value = 2299.2501
$ws.Range("L122").Value = 1700.0001
$ws.Range("M122").Value = 150.7498999999998
$ws.Range("N122").Value = -6600.0001

$ws.Range("H138").Value = 2413.3965
$ws.Range("J138").Value = 2489.3022
$ws.Range("L138").Value = 7467.9066
$ws.Range("N138").Value = -17747.9066

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 9004
$ws.Range("I9").Value = 7999
$ws.Range("K9").Value = 7999
$ws.Range("M9").Value = -7829

$ws.Range("H20").Value = 9004
$ws.Range("I20").Value = 7999
$ws.Range("K20").Value = 7999
$ws.Range("M20").Value = -7729

$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12286.588
$ws.Range("J31").Value = 6921.4
$ws.Range("L31").Value = 6921.4
$ws.Range("N31").Value = -7511.4

$ws.Range("H34").Value = 12286.588
$ws.Range("J34").Value = 6921.4
$ws.Range("L34").Value = 6921.4
$ws.Range("N34").Value = -7325.4

$ws.Range("H58").Value = 11698.479
$ws.Range("J58").Value = 34208.332
$ws.Range("L58").Value = 34208.332
$ws.Range("N58").Value = -34614.332

$ws.Range("H136").Value = 11698.479
$ws.Range("J136").Value = 34208.332
$ws.Range("L136").Value = 102624.996
$ws.Range("N136").Value = -107724.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 360
$ws.Range("I47").Value = 360
$ws.Range("K47").Value = 1080
$ws.Range("M47").Value = -649

$ws.Range("H104").Value = 1840
$ws.Range("I104").Value = 1840
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 5520
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -2899
$ws.Range("N104").ClearContents()

$ws.Range("H113").Value = 5124.0454
$ws.Range("I113").Value = 10578.1
$ws.Range("J113").Value = 579
$ws.Range("K113").Value = 31734.3
$ws.Range("L113").Value = 1737
$ws.Range("M113").Value = -29564.3
$ws.Range("N113").Value = -6077

$ws.Range("H122").Value = 453.79166
$ws.Range("I122").Value = 240.17647
$ws.Range("K122").Value = 2161.58823
$ws.Range("M122").Value = 288.4117700000002

$ws.Range("H131").Value = 787.92
$ws.Range("J131").Value = 787.92
$ws.Range("L131").Value = 2363.76
$ws.Range("N131").Value = -12443.76

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H11").Value = 4599999
$ws.Range("I11").Value = 6000000
$ws.Range("K11").Value = 6000000
$ws.Range("M11").Value = -5999861

$ws.Range("H21").Value = 311111.12
$ws.Range("I21").Value = 311111.12
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 311111.12
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -310938.12
$ws.Range("N21").ClearContents()

$ws.Range("H30").Value = 311111.12
$ws.Range("I30").Value = 311111.12
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 311111.12
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -311006.12
$ws.Range("N30").ClearContents()

$ws.Range("H80").Value = 3722.1365
$ws.Range("I80").Value = 3072.9167
$ws.Range("J80").Value = 4501.2
$ws.Range("K80").Value = 3072.9167
$ws.Range("L80").Value = 4501.2
$ws.Range("M80").Value = -2074.9167
$ws.Range("N80").Value = -6497.2

$ws.Range("H83").Value = 3722.1365
$ws.Range("I83").Value = 3072.9167
$ws.Range("J83").Value = 4501.2
$ws.Range("K83").Value = 15364.5835
$ws.Range("L83").Value = 22506
$ws.Range("M83").Value = -10372.5835
$ws.Range("N83").Value = -32490

$ws.Range("H113").Value = 1764.2
$ws.Range("I113").Value = 1611.1666
$ws.Range("J113").Value = 1993.75
$ws.Range("K113").Value = 1611.1666
$ws.Range("L113").Value = 1993.75
$ws.Range("M113").Value = 558.8334
$ws.Range("N113").Value = -6333.75

$ws.Range("H122").Value = 133336110
$ws.Range("I122").Value = 166667570
$ws.Range("J122").Value = 125003250
$ws.Range("K122").Value = 500002710
$ws.Range("L122").Value = 375009750
$ws.Range("M122").Value = -500000260
$ws.Range("N122").Value = -375014650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4961.12
$ws.Range("I7").Value = 3609.8462
$ws.Range("J7").Value = 6425
$ws.Range("K7").Value = 3609.8462
$ws.Range("L7").Value = 6425
$ws.Range("M7").Value = -3497.8462
$ws.Range("N7").Value = -6649

$ws.Range("H126").Value = 4961.12
$ws.Range("I126").Value = 3609.8462
$ws.Range("J126").Value = 6425
$ws.Range("K126").Value = 10829.5386
$ws.Range("L126").Value = 19275
$ws.Range("M126").Value = -8359.5386
$ws.Range("N126").Value = -24215

$ws.Range("H132").Value = 1930.5862
$ws.Range("I132").Value = 1341.6316
$ws.Range("K132").Value = 4024.8948
$ws.Range("M132").Value = -1494.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 70003
$ws.Range("J2").Value = 70003
$ws.Range("L2").Value = 70003
$ws.Range("N2").Value = -70227
